# CORRIDAS TK MULTIMARCAS.xlsx — update "MES 01" trip log:
#   Row 63 was a blank trip entry (count 0, no neighbourhood). Fill it in
#   with a 15-trip entry for "LAGO AZUL " (the sheet's running total in
#   B66 recalculates automatically from 462 to 477).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")

$ws.Range("A63").Value = 15
$ws.Range("B63").Value = "LAGO AZUL "

# Keep the sheet's active selection/view consistent with where the user
# was working after making the edit.
$ws.Range("D65").Select() | Out-Null

$excel.Calculate() | Out-Null
